# Rename the embedded logo pictures in the document's footers/header.
#
#   footer (id=3 / first-page footer, PearsonLogo)  : image1.png -> image2.png
#   footer (id=2 / default footer,   PearsonLogo)   : image1.png -> image2.png
#   header (id=1 / first-page header, BTec_Logo-Orange): image2.jpg -> image1.jpg
#
# InlineShape has no settable "Name" in the Word object model (only the
# floating Shape object does), so each inline picture is round-tripped
# through ConvertToShape()/ConvertToInlineShape() to reach the Name
# property, then restored to an inline picture in place.

$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-InlinePicture($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

# First-page footer (docPr id="3") - Pearson logo
$footerFirst = $sec.Footers.Item(2)
Rename-InlinePicture $footerFirst.Range.InlineShapes.Item(1) "image2.png"

# Default footer (docPr id="2") - Pearson logo
$footerDefault = $sec.Footers.Item(1)
Rename-InlinePicture $footerDefault.Range.InlineShapes.Item(1) "image2.png"

# First-page header (docPr id="1") - BTec logo
$headerFirst = $sec.Headers.Item(2)
Rename-InlinePicture $headerFirst.Range.InlineShapes.Item(1) "image1.jpg"

Write-Output "Renamed logo inline shapes in footers/header."
